$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.706.63'
$ws.Range("E2").Value = '  -0.83%  '

$ws.Range("D3").Value = '2.297.31'
$ws.Range("E3").Value = '  +2.46%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '95.04'
$ws.Range("E5").Value = '  -3.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.16'
$ws.Range("E6").Value = '  -1.69%  '

$ws.Range("E7").Value = '  -0.82%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("E9").Value = '  -5.73%  '

$ws.Range("E10").Value = '  -9.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0936'
$ws.Range("E11").Value = '  -1.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.77'
$ws.Range("E12").Value = '  -9.21%  '

$ws.Range("E13").Value = '  -0.18%  '

$ws.Range("D14").Value = '2.643.71'
$ws.Range("E14").Value = '  +22.22%  '

$ws.Range("E15").Value = '  -1.25%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.851'
$ws.Range("E16").Value = '  +3.10%  '

$ws.Range("D17").Value = '2.310.75'
$ws.Range("E17").Value = '  +2.72%  '

$ws.Range("D18").Value = '43.671.97'
$ws.Range("E18").Value = '  -0.84%  '

$ws.Range("E19").Value = '  +2.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").Value = '  -0.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.46'
$ws.Range("E21").Value = '  +1.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.37'
$ws.Range("E22").Value = '  +2.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.84'
$ws.Range("E23").Value = '  +0.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.96'
$ws.Range("E24").Value = '  -3.16%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.35'
$ws.Range("E26").Value = '  -2.34%  '

$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.50'
$ws.Range("E27").Value = '  -1.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.46'
$ws.Range("E28").Value = '  -2.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("E29").Value = '  +0.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.72'
$ws.Range("E30").Value = '  +1.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.66'
$ws.Range("E31").Value = '  -7.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.97'
$ws.Range("E32").Value = '  +4.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0885'
$ws.Range("E33").Value = '  -4.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.39'
$ws.Range("E34").Value = '  -3.85%  '

$ws.Range("E35").Value = '  +0.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.108'
$ws.Range("E36").Value = '  -5.14%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.44'
$ws.Range("E37").Value = '  +1.30%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0354'
$ws.Range("E38").Value = '  +0.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.25'
$ws.Range("E39").Value = '  -12.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.36'
$ws.Range("E40").Value = '  +8.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.235'
$ws.Range("E41").Value = '  -6.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.36'
$ws.Range("E42").Value = '  +17.02%  '

$ws.Range("E43").Value = '  -8.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.18'
$ws.Range("E44").Value = '  -1.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.93'
$ws.Range("E45").Value = '  +5.90%  '

$ws.Range("E46").Value = '  -4.08%  '

$ws.Range("E47").Value = '  +0.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.89'
$ws.Range("E48").Value = '  -2.83%  '

$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.54'
$ws.Range("E50").Value = '  +7.53%  '

$ws.Range("D51").Value = '2.521.74'
$ws.Range("E51").Value = '  +3.16%  '
